# ValueSet-fr-route-of-administration.xlsx update
#
#  Metadata sheet:
#   - URL: hl7.fr/fhir/fr/medication -> hl7.fr/ig/fhir/medication
#   - Date: 2025-04-10T15:35:36+00:00 -> 2026-01-15T08:54:26+00:00
#   - Jurisdiction: (empty) -> FRANCE
#   - Description: "FrMedicationRequest" -> "FRMedicationRequest"
#
#  Include #0 sheet: the old "Property/Operation/Value" constraint table
#  (3 columns) is replaced by a 2-column "Codes" / "System URI" summary
#  table.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "Metadata"
$ws2 = $wb.Worksheets.Item(2)   # "Include #0"

# --- Metadata sheet updates -------------------------------------------------

$ws1.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/ValueSet/fr-route-of-administration"
$ws1.Range("B8").Value = "2026-01-15T08:54:26+00:00"
$ws1.Range("B11").Value = "FRANCE"
$ws1.Range("B12").Value = "Le jeu de valeurs à utiliser pour coder l'élément *dosageInstruction.route* de la ressource *FRMedicationRequest*."

# --- Include #0 sheet: rebuild the table ------------------------------------

# Wipe the old 3-column x 4-row table entirely (clean slate - no leftover
# cells/columns), then lay out only the cells the new 2-column table needs.
$ws2.Range("A1:C4").EntireRow.Delete()

$ws2.Range("A1").Value = "Codes"
$ws2.Range("A2").Value = "All codes"
$ws2.Range("A3").Value = ""
$ws2.Range("B3").Value = ""
$ws2.Range("A4").Value = "System URI"
$ws2.Range("B4").Value = "http://standardterms.edqm.eu"

# Re-apply the sheet's header / body cell styles (bold header row, wrapped
# body rows) by copying formats from the Metadata sheet, which uses the
# same two cell styles (s="1" header, s="2" body).
$ws1.Range("A1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)
$ws2.Range("A3:B3").PasteSpecial(-4122)
$ws2.Range("A4:B4").PasteSpecial(-4122)

$excel.CutCopyMode = 0
